$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# A new incident (Caso -627, "Av San Martin 3231") was reported and now sits
# at sheet row 87, ahead of the previous row 87 (Caso 7371). Insert a fresh
# row there — this shifts every row from 87..93 down to 88..94, exactly like
# the diff shows (dimension grows from A1:R93 to A1:R94).
$ws.Rows("87:87").Insert()

# Text columns must stay text even though several look numeric (Caso, OT,
# Comuna, etc. are all stored as text in this sheet) — prefix with an
# apostrophe so Excel keeps them as strings instead of coercing to numbers.
$ws.Range("A87").Value = "'-627"
$ws.Range("B87").Value = "'10/1/2025"
$ws.Range("C87").Value = "Av San Martin 3231"
$ws.Range("D87").Value = "'15"
$ws.Range("E87").Value = "'810093647"
$ws.Range("F87").Value = "PEBCOM"
$ws.Range("G87").Value = "Pendiente"
$ws.Range("H87").Value = "Columna chocada"
$ws.Range("I87").Value = 1
$ws.Range("J87").Value = "Cambio"
$ws.Range("K87").Value = "Sin equipos"
$ws.Range("L87").Value = "Pasante"
$ws.Range("M87").Value = -58.469321
$ws.Range("N87").Value = -34.601663
$ws.Range("O87").Value = "Paternal"
$ws.Range("P87").Value = "Capital Norte"
$ws.Range("Q87").Value = "NRA-F"
$ws.Range("R87").Value = "ARATO-25058.PO.1NRA"
